$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.065.04'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.324.50'
$ws.Range('E3').Value = '  +0.89%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.997'
$ws.Range('E4').Value = '  -0.38%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.16'
$ws.Range('E5').Value = '  +4.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '185.15'
$ws.Range('E6').Value = '  -1.92%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.319.21'
$ws.Range('E8').Value = '  +1.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.575'
$ws.Range('E9').Value = '  -2.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.180'
$ws.Range('E10').Value = '  -2.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.577'
$ws.Range('E11').Value = '  -1.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '46.94'
$ws.Range('E12').Value = '  -1.30%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000267'
$ws.Range('E13').Value = '  -0.85%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '664.33'
$ws.Range('E14').Value = '  +9.87%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.850.47'
$ws.Range('E15').Value = '  +0.66%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.47'
$ws.Range('E16').Value = '  -1.96%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '66.147.66'
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.90'
$ws.Range('E18').Value = '  -0.40%  '
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.117'
$ws.Range('E19').Value = '  -0.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.321.89'
$ws.Range('E20').Value = '  +0.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.10'
$ws.Range('E21').Value = '  +0.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.896'
$ws.Range('E22').Value = '  -1.53%  '
$ws.Range('E23').Value = '  -2.59%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '101.57'
$ws.Range('E24').Value = '  +0.71%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.04'
$ws.Range('E25').Value = '  -0.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.98'
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.78'
$ws.Range('E27').Value = '  +1.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.52'
$ws.Range('E28').Value = '  -1.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '31.54'
$ws.Range('E29').Value = '  +4.47%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.45'
$ws.Range('E30').Value = '  -1.92%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.71'
$ws.Range('E31').Value = '  +0.19%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '592.12'
$ws.Range('E32').Value = '  +5.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.81'
$ws.Range('E33').Value = '  -5.67%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '10.99'
$ws.Range('E34').Value = '  -0.71%  '
$ws.Range('E35').Value = '  +0.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.841.74'
$ws.Range('E36').Value = '  +3.82%  '
$ws.Range('E37').Value = '  +0.25%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '55.92'
$ws.Range('E38').Value = '  -2.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.68'
$ws.Range('E39').Value = '  -0.76%  '
$ws.Range('B40').Value = 'PEPE'
$ws.Range('C40').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0698'
$ws.Range('E40').Value = '  -3.58%  '
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '32.76'
$ws.Range('E41').Value = '  -3.18%  '
$ws.Range('E42').Value = '  -3.11%  '
$ws.Range('B43').Value = 'ApeXProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.43'
$ws.Range('E43').Value = '  +5.80%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.16'
$ws.Range('E44').Value = '  -4.09%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.335'
$ws.Range('E45').Value = '  -1.37%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0412'
$ws.Range('E46').Value = '  -2.71%  '
$ws.Range('E47').Value = '  -10.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.128'
$ws.Range('E48').Value = '  -1.43%  '
$ws.Range('E49').Value = '  +0.21%  '
$ws.Range('E50').Value = '  -1.78%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.30'
$ws.Range('E51').Value = '  +1.90%  '
